# "book texture & mapping 1"
# Adds a second table below the existing ghost-types table, listing
# investigators (Stejskal, Skarka, Lazarov, Bobek) as column headers in
# row 11, plus a few more evidence/attribute rows (12-15) with "x" marks,
# mirroring the layout of the first table (rows 3-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 11) - investigator names.
# Entry order below reproduces the original authoring order so that the
# shared-string table ends up built in the same sequence as the source
# workbook (B11, A12, A13, then back to C11/D11/E11, then A14, A15).
$ws.Range("B11").Value = "Stejskal"
$ws.Range("A12").Value = "Bakalari"
$ws.Range("A13").Value = "Code writing"
$ws.Range("C11").Value = "Skarka"
$ws.Range("D11").Value = "Lazarov"
$ws.Range("E11").Value = "Bobek"
$ws.Range("A14").Value = "Virtualka"
$ws.Range("A15").Value = "Speech"

# Row 12 - Bakalari: x under Skarka, Lazarov, Bobek
$ws.Range("C12").Value = "x"
$ws.Range("D12").Value = "x"
$ws.Range("E12").Value = "x"

# Row 13 - Code writing: x under Stejskal
$ws.Range("B13").Value = "x"

# Row 14 - Virtualka: x under Skarka, Lazarov
$ws.Range("C14").Value = "x"
$ws.Range("D14").Value = "x"

# Row 15 - Speech: x under Stejskal, Lazarov, Bobek
$ws.Range("B15").Value = "x"
$ws.Range("D15").Value = "x"
$ws.Range("E15").Value = "x"

# Match the selection left behind in the saved workbook.
$ws.Range("F13").Select()
